$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cinema")

# Insert a new column before column S (19) - shifts S:AO (and beyond) right by one.
$ws.Columns("S:S").Insert()

# Header for new column S, matching style of neighboring header R1 (s="4")
$ws.Range("S1").Value = "Sub brand"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The AutoFilter range does not auto-grow when a column is inserted inside
# its span, so refresh it to cover the new column (A1:AP54).
$ws.AutoFilterMode = $false
$ws.Range("A1:AP54").AutoFilter()

# The backing defined names for the filter database also need to be
# refreshed explicitly to the new range.
$newRef = "=Cinema!`$A`$1:`$AP`$54"
$wb.Names.Item("Cinema!_FilterDatabase").RefersTo = $newRef
$wb.Names.Item("Cinema!_FilterDatabase_0").RefersTo = $newRef

# Update selection to match target state
$ws.Range("S2").Select()
